$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at position 13 ---
# This shifts the old rows 13-23 ("Programa resumido:" ... the requisitos row) down
# to rows 14-24, matching the target A1:C24 layout (was A1:C23).
$ws.Rows.Item(13).Insert()

# The inserted row leaves a phantom styled-but-empty cell in column A; the target
# row 13 has no A cell at all (only B13/C13 hold the professor name), so clear it.
$ws.Range("A13").Clear()

# Give B13/C13 the same look as the other data cells in columns B/C (regular/red
# wrapped text) by copying formats from row 10, then write the actual value.
$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C13").PasteSpecial(-4122)

# --- Row 10 (Objetivos:): replace the misplaced professor name with the real objectives text ---
$ws.Range("B10").Value = 'Promover o conhecimento da Bioquímica no que tange à bioenergética e metabolismoCompreender a importância das macromoléculas e de suas vias metabólicasUtilizar os conhecimentos como pré-requisito para as disciplinas do curso de Engenharia Bioquímica'
$ws.Range("C10").Value = 'Promover o conhecimento da Bioquímica no que tange à bioenergética e metabolismoCompreender a importância das macromoléculas e de suas vias metabólicasUtilizar os conhecimentos como pré-requisito para as disciplinas do curso de Engenharia Bioquímica'

# --- Row 13 (new, under Docentes responsáveis:): the professor name now lives here ---
$ws.Range("B13").Value = '5840494 - Maria Eleonora Andrade de Carvalho'
$ws.Range("C13").Value = '5840494 - Maria Eleonora Andrade de Carvalho'

# --- Row 14 (Programa resumido:): replace "Semestral" with the real short syllabus (PT) ---
$ws.Range("B14").Value = 'Introdução ao metabolismo, BioenergéticaOxidações biológicas, TransporteGlicídios - metabolismoFotossínteseLipídeos - metabolismoAminoácidos - metabolismoIntegração MetabólicaCiclos vitais: oxigênio, carbono, nitrogênio e enxofre'
$ws.Range("C14").Value = 'Introdução ao metabolismo, BioenergéticaOxidações biológicas, TransporteGlicídios - metabolismoFotossínteseLipídeos - metabolismoAminoácidos - metabolismoIntegração MetabólicaCiclos vitais: oxigênio, carbono, nitrogênio e enxofre'

# --- Row 16 (Programa:): replace the misplaced date with the real full syllabus (PT) ---
$ws.Range("B16").Value = 'Introdução ao metabolismo. Proteínas, polissacarídeos, lipídios: vias catabólicas e anabólicas. Bioenergética. Variação de energia livre: relação com a constante de equilíbrio e com o potencial redox. Processos exergônicos. Papel do fosfato: potencial de transferência de grupo fosfato. Importância energética do ATP.Oxidações biológicas. Coenzimas transportadoras de prótons e elétrons: nucleotídeos, flavino nucleotídeos, coenzima Q. Desidrogenase piridino e flavino nucleotídeos dependentes. Oxidases. Estrutura da membrana mitocondrial. Cadeia respiratória: função. Fosforilação oxidativa.Transporte. Composição das membranas biológicas: constituição química, caráter ?barreira permeabilidade?. Carreadores e canais, ionóforos. Transporte: mediado e não mediado. Glicídios - metabolismo. Degradação anaeróbica e aeróbica de glicídios: glicólise - localização das enzimas operantes, reações, produção de NADH.H+, fosforilação ao nível de substrato, balanço energético; ciclo de Krebs - localização das enzimas operantes, reações, produção de coenzimas reduzidas, balanço energético. ?Shunt? das hexoses-fosfato (ciclo das pentoses): localização das enzimas operantes, reações (fases oxidativa e não oxidativa), produção de NADPH.H  (implicação fisiológica). Fermentações: definição, fermentação e respiração, matérias primas usadas em fermentação amilácea e sacarínea, agente de fermentação, fermentações anaeróbicas - alcoólica e lática, fermentações aeróbicas - acética e cítrica. Fotossíntese. Estrutura dos cloroplastos. Luz: energia eletromagnética. Papel da clorofila na fotossíntese. Fotofosforilação cíclica e não cíclica. Redução do NADP . Fotólise da água. Síntese do aceptor de CO2, Ru-1,5diP. Ciclo de Calvin.Lipídios - metabolismo. b-oxidação de ácidos graxos de cadeia: com número par de átomos de C, com número impar de átomos de C, ramificada, balanço energético da b-oxidação. a-oxidação e w-oxidação. Metabolismo do glicerol. Formação de corpos cetônicos. Biossíntese de ácidos graxos.Aminoácidos - vias catabólicas. Digestão de proteínas, enzimas envolvidas e zimogênios, absorção. Transaminação, desaminação oxidativa, aminoácidos cetogênicos e glicogênicos, descarboxilação, ciclo da uréia. Eliminação de nitrogênio, vertebrados ureotélicos.Integração metabólica. Interelação do metabolismo intermediário de glicídios, lipídios, aminoácidos e ácidos nucléicos. Metabólitos comuns ao metabolismo de glicídios, lipídios e aminoácidos.Ciclos vitais: oxigênio, carbono, nitrogênio e enxofre. Ciclos do oxigênio e do carbono. Ciclo do nitrogênio: fixação biológica, nitrificação, utilização do nitrato, incorporação de amônia em compostos orgânicos. Ciclo do enxofre: assimilação do sulfato'
$ws.Range("C16").Value = 'Introdução ao metabolismo. Proteínas, polissacarídeos, lipídios: vias catabólicas e anabólicas. Bioenergética. Variação de energia livre: relação com a constante de equilíbrio e com o potencial redox. Processos exergônicos. Papel do fosfato: potencial de transferência de grupo fosfato. Importância energética do ATP.Oxidações biológicas. Coenzimas transportadoras de prótons e elétrons: nucleotídeos, flavino nucleotídeos, coenzima Q. Desidrogenase piridino e flavino nucleotídeos dependentes. Oxidases. Estrutura da membrana mitocondrial. Cadeia respiratória: função. Fosforilação oxidativa.Transporte. Composição das membranas biológicas: constituição química, caráter ?barreira permeabilidade?. Carreadores e canais, ionóforos. Transporte: mediado e não mediado. Glicídios - metabolismo. Degradação anaeróbica e aeróbica de glicídios: glicólise - localização das enzimas operantes, reações, produção de NADH.H+, fosforilação ao nível de substrato, balanço energético; ciclo de Krebs - localização das enzimas operantes, reações, produção de coenzimas reduzidas, balanço energético. ?Shunt? das hexoses-fosfato (ciclo das pentoses): localização das enzimas operantes, reações (fases oxidativa e não oxidativa), produção de NADPH.H  (implicação fisiológica). Fermentações: definição, fermentação e respiração, matérias primas usadas em fermentação amilácea e sacarínea, agente de fermentação, fermentações anaeróbicas - alcoólica e lática, fermentações aeróbicas - acética e cítrica. Fotossíntese. Estrutura dos cloroplastos. Luz: energia eletromagnética. Papel da clorofila na fotossíntese. Fotofosforilação cíclica e não cíclica. Redução do NADP . Fotólise da água. Síntese do aceptor de CO2, Ru-1,5diP. Ciclo de Calvin.Lipídios - metabolismo. b-oxidação de ácidos graxos de cadeia: com número par de átomos de C, com número impar de átomos de C, ramificada, balanço energético da b-oxidação. a-oxidação e w-oxidação. Metabolismo do glicerol. Formação de corpos cetônicos. Biossíntese de ácidos graxos.Aminoácidos - vias catabólicas. Digestão de proteínas, enzimas envolvidas e zimogênios, absorção. Transaminação, desaminação oxidativa, aminoácidos cetogênicos e glicogênicos, descarboxilação, ciclo da uréia. Eliminação de nitrogênio, vertebrados ureotélicos.Integração metabólica. Interelação do metabolismo intermediário de glicídios, lipídios, aminoácidos e ácidos nucléicos. Metabólitos comuns ao metabolismo de glicídios, lipídios e aminoácidos.Ciclos vitais: oxigênio, carbono, nitrogênio e enxofre. Ciclos do oxigênio e do carbono. Ciclo do nitrogênio: fixação biológica, nitrificação, utilização do nitrato, incorporação de amônia em compostos orgânicos. Ciclo do enxofre: assimilação do sulfato'

# --- Row 19 (Método:): replace the misplaced professor name with the real method text ---
$ws.Range("B19").Value = 'A avaliação será feita por meio de provas escritas.'
$ws.Range("C19").Value = 'A avaliação será feita por meio de provas escritas.'

# --- Row 20 (Critério:): the final-grade formula text belongs here ---
$ws.Range("B20").Value = 'A Nota final (NF) será calculada da seguinte maneira:NF = (P1 + P2)/2'
$ws.Range("C20").Value = 'A Nota final (NF) será calculada da seguinte maneira:NF = (P1 + P2)/2'

# --- Row 21 (Norma de recuperação:): the make-up exam text belongs here ---
$ws.Range("B21").Value = 'A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2'
$ws.Range("C21").Value = 'A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2'

# --- Row 22 (Bibliografia:): replace the misplaced text with the real bibliography ---
$ws.Range("B22").Value = '1.Nelson, D.L.; Cox, M.M. Lehninger Principles of Biochemistry. Third Edition, Worth Publisher, New York, 20002.Voet, D; Voet, J G.; Pratt, C.W. Fundamentos de Bioquímica. Editora ARTMED, Porto Alegre, 20003.Stryer, L. Biochemistry. W.H. Freeman Company, New York, 19884.Jain, M.K. Introduction to Biological Membranes. John Wiley & Sons Inc., New York, 1988'
$ws.Range("C22").Value = '1.Nelson, D.L.; Cox, M.M. Lehninger Principles of Biochemistry. Third Edition, Worth Publisher, New York, 20002.Voet, D; Voet, J G.; Pratt, C.W. Fundamentos de Bioquímica. Editora ARTMED, Porto Alegre, 20003.Stryer, L. Biochemistry. W.H. Freeman Company, New York, 19884.Jain, M.K. Introduction to Biological Membranes. John Wiley & Sons Inc., New York, 1988'
